$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Table 1 "All Major Components" (rows 7-10)
# ---------------------------------------------------------------------------

# Row 9 - "Green LED": part number, qty and current draw change; G9 becomes a
# live formula (E9*F9) instead of a hard-coded value.
$ws.Range("C9").Value = "150080VS75000"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 100

$ws.Range("G8").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Formula = "=E9*F9"

# Row 10 - brand new "Red LED" row (was completely blank). Clone the
# formatting of row 9 first, then fix up the couple of cells that use a
# different style, then fill in the values.
$ws.Range("A9:H9").Copy()
$ws.Range("A10:H10").PasteSpecial(-4122)

$ws.Range("C8").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("G8").Copy()
$ws.Range("G10").PasteSpecial(-4122)

$ws.Range("B10").Value = "Red LED"
$ws.Range("C10").Value = "150080RS75000"
$ws.Range("D10").Value = " +3.3V"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 100
$ws.Range("G10").Formula = "=E10*F10"
$ws.Range("H10").Value = "mA"

# ---------------------------------------------------------------------------
# Table 2 "+3.3V Power Rail" (rows 11-18)
# ---------------------------------------------------------------------------

# Row 13 - "Green LED": same part-number / qty / current changes as row 9.
# (C13 shares the same underlying string as C9, already updated above, but
# set it explicitly too in case the two cells aren't literally linked.)
$ws.Range("C13").Value = "150080VS75000"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 100

# Row 14 - was blank except for G14/H14; fill in the new "Red LED" entry.
$ws.Range("B13:F13").Copy()
$ws.Range("B14:F14").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("B14").Value = "Red LED"
$ws.Range("C14").Value = "150080RS75000"
$ws.Range("D14").Value = " +3.3V"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 100

# ---------------------------------------------------------------------------
# External power source section
# ---------------------------------------------------------------------------

# Row 24 - Power Source 1 Selection: new supply part number + higher current.
$ws.Range("C24").Value = "ALT-1208"
$ws.Range("F24").Value = 8000

# Row 26 - Power Rails Connected to External Power Source 1: new regulator.
$ws.Range("C26").Value = "AP62300TWU-7"

# ---------------------------------------------------------------------------
# Cosmetic: A29 picks up a plain white/Arial style in the source edit.
# ---------------------------------------------------------------------------
$ws.Range("A29").Style = "Normal"
$ws.Range("A29").Font.Name = "Arial"
$ws.Range("A29").Font.Size = 11
$ws.Range("A29").Font.Bold = $false
$ws.Range("A29").Font.Color = 3355443
$ws.Range("A29").Interior.Pattern = 1
$ws.Range("A29").Interior.Color = 16777215

$wb.Application.CalculateFull()
